$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = -8.615999999999987
$ws.Range("C12").Value = -14.62430000000001
$ws.Range("D12").Value = -8.092600000000003
$ws.Range("D14").Value = -8.707199999999998
$ws.Range("D22").Value = -7.856099999999997
